# CryCompanywiseStockReport_1.xlsx - stock quantity/value corrections.
#
# Each product line has: D = unit cost, E = unit price, F = quantity,
# G = line value (= D * F). Updating a line's quantity/value also ripples
# into that company's "Sub Total:" row (B = SUM of the group's G values),
# the overall "Sub Total:" row (B724) and the final "Grand Total:" row
# (B725), which are plain numeric cells (not live formulas) in this sheet -
# so every affected cell is written explicitly below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 29, AL-dazller kreamy kajal (qty 2 -> 1) + company Sub Total (row 40) ---
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 40.99
$ws.Range("B40").Value = 50643.81

# --- row 47 (qty 170 -> 168) ---
$ws.Range("F47").Value = 168
$ws.Range("G47").Value = 32405.52

# --- row 64 (qty 41 -> 39) ---
$ws.Range("F64").Value = 39
$ws.Range("G64").Value = 3039.27

# --- row 66 (qty 1 -> 0) ---
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0

# --- row 67 (qty 176 -> 174) + company Sub Total (row 72) ---
$ws.Range("F67").Value = 174
$ws.Range("G67").Value = 45367.02
$ws.Range("B72").Value = 159843.52

# --- row 120 (qty 13 -> 9) + company Sub Total (row 129) ---
$ws.Range("F120").Value = 9
$ws.Range("G120").Value = 420.66
$ws.Range("B129").Value = 65467.56

# --- row 181 (qty 22 -> 21) + company Sub Total (row 199) ---
$ws.Range("F181").Value = 21
$ws.Range("G181").Value = 6090.21
$ws.Range("B199").Value = 52189.5

# --- row 216 (qty 95 -> 94) ---
$ws.Range("F216").Value = 94
$ws.Range("G216").Value = 5114.54

# --- row 219 (qty 187 -> 185) ---
$ws.Range("F219").Value = 185
$ws.Range("G219").Value = 23435.8

# --- row 221 (qty 139 -> 137) + company Sub Total (row 224) ---
$ws.Range("F221").Value = 137
$ws.Range("G221").Value = 15383.73
$ws.Range("B224").Value = 60374.52

# --- row 270 (qty 46 -> 42) ---
$ws.Range("F270").Value = 42
$ws.Range("G270").Value = 1463.28

# --- row 273 (qty 97 -> 95) + company Sub Total (row 301) ---
$ws.Range("F273").Value = 95
$ws.Range("G273").Value = 4035.6
$ws.Range("B301").Value = 90431.09

# --- rows 303 / 304: the two lots' code/price/qty/value are swapped ---
$ws.Range("B303").Value = 63565
$ws.Range("E303").Value = 109.19
$ws.Range("F303").Value = 60
$ws.Range("G303").Value = 6162.6

$ws.Range("B304").Value = 61610
$ws.Range("E304").Value = 122.71
$ws.Range("F304").Value = -58
$ws.Range("G304").Value = -5957.18

# --- row 308 (qty 34 -> 33) ---
$ws.Range("F308").Value = 33
$ws.Range("G308").Value = 3769.59

# --- rows 312 / 313: the two lots' code/price/qty/value are swapped ---
$ws.Range("B312").Value = 57802
$ws.Range("E312").Value = 162.71
$ws.Range("F312").Value = -79
$ws.Range("G312").Value = -11334.92

$ws.Range("B313").Value = 63531
$ws.Range("E313").Value = 152.53
$ws.Range("F313").Value = 20
$ws.Range("G313").Value = 2869.6

# --- company Sub Total (row 334) ---
$ws.Range("B334").Value = -24176.41

# --- row 357 (qty 171 -> 170) + company Sub Total (row 362) ---
$ws.Range("F357").Value = 170
$ws.Range("G357").Value = 24588.8
$ws.Range("B362").Value = 66894.08

# --- row 376 (qty 148 -> 146) + company Sub Total (row 378) ---
$ws.Range("F376").Value = 146
$ws.Range("G376").Value = 24234.54
$ws.Range("B378").Value = 44091.02

# --- row 414 (qty 151 -> 150) + company Sub Total (row 423) ---
$ws.Range("F414").Value = 150
$ws.Range("G414").Value = 2377.5
$ws.Range("B423").Value = 151843.36

# --- row 436 (qty 183 -> 181) ---
$ws.Range("F436").Value = 181
$ws.Range("G436").Value = 8376.68

# --- row 438 (qty 43 -> 42) + company Sub Total (row 444) ---
$ws.Range("F438").Value = 42
$ws.Range("G438").Value = 2033.22
$ws.Range("B444").Value = 18520.06

# --- row 461 (qty 29 -> 28) + company Sub Total (row 464) ---
$ws.Range("F461").Value = 28
$ws.Range("G461").Value = 6222.44
$ws.Range("B464").Value = 76839.3

# --- rows 485 / 486: the two lots' code/price/qty/value are swapped ---
$ws.Range("B485").Value = 64810
$ws.Range("E485").Value = 291.22
$ws.Range("F485").Value = 0
$ws.Range("G485").Value = 0

$ws.Range("B486").Value = 53319
$ws.Range("E486").Value = 310.64
$ws.Range("F486").Value = -6
$ws.Range("G486").Value = -1643.52

# --- rows 512 / 513: the two lots' code/price/qty/value are swapped ---
$ws.Range("B512").Value = 64830
$ws.Range("E512").Value = 34.9
$ws.Range("F512").Value = 83
$ws.Range("G512").Value = 2724.89

$ws.Range("B513").Value = 60022
$ws.Range("E513").Value = 37.22
$ws.Range("F513").Value = -113
$ws.Range("G513").Value = -3709.79

# --- row 517 (qty 165 -> 164) ---
$ws.Range("F517").Value = 164
$ws.Range("G517").Value = 16378.68

# --- row 518 (qty 7 -> 6) + company Sub Total (row 531) ---
$ws.Range("F518").Value = 6
$ws.Range("G518").Value = 711.48
$ws.Range("B531").Value = 104104.18

# --- row 538 (qty 2 -> 1) + company Sub Total (row 541) ---
$ws.Range("F538").Value = 1
$ws.Range("G538").Value = 43.18
$ws.Range("B541").Value = 16209.26

# --- row 564 (qty 110 -> 109) + company Sub Total (row 567) ---
$ws.Range("F564").Value = 109
$ws.Range("G564").Value = 13281.65
$ws.Range("B567").Value = 14907.27

# --- row 645 (qty 1 -> 0) ---
$ws.Range("F645").Value = 0
$ws.Range("G645").Value = 0

# --- company Sub Total (row 663) ---
$ws.Range("B663").Value = 58482.24

# --- row 680 (qty 252 -> 233) + company Sub Total (row 686) ---
$ws.Range("F680").Value = 233
$ws.Range("G680").Value = 38004.63
$ws.Range("B686").Value = 39017.18

# --- row 712 (qty 102 -> 101) + company Sub Total (row 719) ---
$ws.Range("F712").Value = 101
$ws.Range("G712").Value = 3958.19
$ws.Range("B719").Value = 53161.1

# --- overall Sub Total (row 724) and Grand Total (row 725) ---
$ws.Range("B724").Value = 2038968.53
$ws.Range("B725").Value = 2038968.53
